# Weekly price-sheet update: a new observation is reported for row 42
# (same item/quality/price-point as the existing row 42, but at a later
# date), which pushes the rest of the table down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 42:116 down to 43:117, carrying formats with them, and
# leave a blank (but correctly formatted) row 42 to populate.
$ws.Rows("42:42").Insert()

$ws.Range("A42").Value = 7
$ws.Range("B42").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C42").Value = 'Ñuble'
$ws.Range("D42").Value = 44994
$ws.Range("E42").Value = 16
$ws.Range("F42").Value = 100112037
$ws.Range("G42").Value = 'Cebollín'
$ws.Range("H42").Value = 'Sin especificar'
$ws.Range("I42").Value = 'Primera'
$ws.Range("J42").Value = 50
$ws.Range("K42").Value = 6000
$ws.Range("L42").Value = 6000
$ws.Range("M42").Value = 6000
$ws.Range("N42").Value = '$/paquete 36 unidades'
$ws.Range("O42").Value = 'Provincia de Diguillín'
$ws.Range("P42").Value = 167
$ws.Range("Q42").Value = 36
$ws.Range("R42").Value = 'Hortaliza'
